$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = ' Sharjah'
$ws.Range("B9").Value = ' October 26 2020'
$ws.Range("C9").Value = 'Kings XI won by 8 wickets (with 7 balls remaining)'
$ws.Range("D9").Value = 'Kolkata Knight Riders'
$ws.Range("E9").Value = 'Kings XI Punjab'
$ws.Range("F9").Value = 'Dinesh Karthik †'
$ws.Range("G9").Value = '''0'
$ws.Range("H9").Value = '''2'
$ws.Range("I9").Value = '''0'
$ws.Range("J9").Value = '''0'
$ws.Range("K9").Value = '''0.00'

# Row 10
$ws.Range("A10").Value = ' Abu Dhabi'
$ws.Range("B10").Value = ' October 18 2020'
$ws.Range("C10").Value = 'Match tied (KKR won the one-over eliminator)'
$ws.Range("D10").Value = 'Kolkata Knight Riders'
$ws.Range("E10").Value = 'Sunrisers Hyderabad'
$ws.Range("F10").Value = 'Dinesh Karthik †'
$ws.Range("G10").Value = '''29'
$ws.Range("H10").Value = '''14'
$ws.Range("I10").Value = '''2'
$ws.Range("J10").Value = '''2'
$ws.Range("K10").Value = '''207.14'

# Row 11
$ws.Range("A11").Value = ' Dubai (DSC)'
$ws.Range("B11").Value = ' October 29 2020'
$ws.Range("C11").Value = 'Super Kings won by 6 wickets'
$ws.Range("D11").Value = 'Kolkata Knight Riders'
$ws.Range("E11").Value = 'Chennai Super Kings'
$ws.Range("F11").Value = 'Dinesh Karthik †'
$ws.Range("G11").Value = '''21'
$ws.Range("H11").Value = '''10'
$ws.Range("I11").Value = '''3'
$ws.Range("J11").Value = '''0'
$ws.Range("K11").Value = '''210.00'

# Row 12
$ws.Range("A12").Value = ' Abu Dhabi'
$ws.Range("B12").Value = ' October 24 2020'
$ws.Range("C12").Value = 'KKR won by 59 runs'
$ws.Range("D12").Value = 'Kolkata Knight Riders'
$ws.Range("E12").Value = 'Delhi Capitals'
$ws.Range("F12").Value = 'Dinesh Karthik †'
$ws.Range("G12").Value = '''3'
$ws.Range("H12").Value = '''6'
$ws.Range("I12").Value = '''0'
$ws.Range("J12").Value = '''0'
$ws.Range("K12").Value = '''50.00'

# Row 13
$ws.Range("A13").Value = ' Abu Dhabi'
$ws.Range("B13").Value = ' October 16 2020'
$ws.Range("C13").Value = 'Mumbai won by 8 wickets (with 19 balls remaining)'
$ws.Range("D13").Value = 'Kolkata Knight Riders'
$ws.Range("E13").Value = 'Mumbai Indians'
$ws.Range("F13").Value = 'Dinesh Karthik †'
$ws.Range("G13").Value = '''4'
$ws.Range("H13").Value = '''8'
$ws.Range("I13").Value = '''1'
$ws.Range("J13").Value = '''0'
$ws.Range("K13").Value = '''50.00'

# Row 14
$ws.Range("A14").Value = ' Dubai (DSC)'
$ws.Range("B14").Value = ' November 01 2020'
$ws.Range("C14").Value = 'KKR won by 60 runs'
$ws.Range("D14").Value = 'Kolkata Knight Riders'
$ws.Range("E14").Value = 'Rajasthan Royals'
$ws.Range("F14").Value = 'Dinesh Karthik †'
$ws.Range("G14").Value = '''0'
$ws.Range("H14").Value = '''1'
$ws.Range("I14").Value = '''0'
$ws.Range("J14").Value = '''0'
$ws.Range("K14").Value = '''0.00'

# Row 15
$ws.Range("A15").Value = ' Abu Dhabi'
$ws.Range("B15").Value = ' October 21 2020'
$ws.Range("C15").Value = 'RCB won by 8 wickets (with 39 balls remaining)'
$ws.Range("D15").Value = 'Kolkata Knight Riders'
$ws.Range("E15").Value = 'Royal Challengers Bangalore'
$ws.Range("F15").Value = 'Dinesh Karthik †'
$ws.Range("G15").Value = '''4'
$ws.Range("H15").Value = '''14'
$ws.Range("I15").Value = '''0'
$ws.Range("J15").Value = '''0'
$ws.Range("K15").Value = '''28.57'

